$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 6.177312000000001
$ws.Range("H2").Value = 18.531936
$ws.Range("I2").Value = 0.1004331027219036
$ws.Range("J2").Value = 0.1004331027219036
$ws.Range("M2").Value = 0.191891
$ws.Range("N2").Value = 0.575673
$ws.Range("O2").Value = 0.07166434268441052
$ws.Range("P2").Value = 0.07166434268441052
$ws.Range("Q2").Value = 1.185370576992
$ws.Range("R2").Value = 10.668335192928
$ws.Range("S2").Value = 0.007197472290321106
$ws.Range("T2").Value = 0.007197472290321106

$ws.Range("G3").Value = 6.177312000000001
$ws.Range("H3").Value = 18.531936
$ws.Range("I3").Value = 0.1004331027219036
$ws.Range("J3").Value = 0.1004331027219036
$ws.Range("O3").Value = 0.4192095090855653
$ws.Range("P3").Value = 0.4192095090855652
$ws.Range("Q3").Value = 6.933973005984001
$ws.Range("R3").Value = 62.40575705385601
$ws.Range("S3").Value = 0.04210251168798938
$ws.Range("T3").Value = 0.04210251168798937

$ws.Range("G4").Value = 6.177312000000001
$ws.Range("H4").Value = 18.531936
$ws.Range("I4").Value = 0.1004331027219036
$ws.Range("J4").Value = 0.1004331027219036
$ws.Range("M4").Value = 1.363254333333333
$ws.Range("N4").Value = 4.089763
$ws.Range("O4").Value = 0.5091261482300243
$ws.Range("P4").Value = 0.5091261482300243
$ws.Range("Q4").Value = 8.421247352351999
$ws.Range("R4").Value = 75.79122617116801
$ws.Range("S4").Value = 0.05113311874359318
$ws.Range("T4").Value = 0.05113311874359318

$ws.Range("I5").Value = 0.7916733321519631
$ws.Range("J5").Value = 0.7916733321519631
$ws.Range("M5").Value = 0.191891
$ws.Range("N5").Value = 0.575673
$ws.Range("O5").Value = 0.07166434268441052
$ws.Range("P5").Value = 0.07166434268441052
$ws.Range("Q5").Value = 9.343794516840001
$ws.Range("R5").Value = 84.09415065156
$ws.Range("S5").Value = 0.05673474896944744
$ws.Range("T5").Value = 0.05673474896944744

$ws.Range("I6").Value = 0.7916733321519631
$ws.Range("J6").Value = 0.7916733321519631
$ws.Range("O6").Value = 0.4192095090855653
$ws.Range("P6").Value = 0.4192095090855652
$ws.Range("S6").Value = 0.3318769889275581
$ws.Range("T6").Value = 0.3318769889275581

$ws.Range("I7").Value = 0.7916733321519631
$ws.Range("J7").Value = 0.7916733321519631
$ws.Range("M7").Value = 1.363254333333333
$ws.Range("N7").Value = 4.089763
$ws.Range("O7").Value = 0.5091261482300243
$ws.Range("P7").Value = 0.5091261482300243
$ws.Range("Q7").Value = 66.38127043403999
$ws.Range("R7").Value = 597.43143390636
$ws.Range("S7").Value = 0.4030615942549577
$ws.Range("T7").Value = 0.4030615942549577

$ws.Range("G8").Value = 6.636180666666667
$ws.Range("H8").Value = 19.908542
$ws.Range("I8").Value = 0.1078935651261332
$ws.Range("J8").Value = 0.1078935651261332
$ws.Range("M8").Value = 0.191891
$ws.Range("N8").Value = 0.575673
$ws.Range("O8").Value = 0.07166434268441052
$ws.Range("P8").Value = 0.07166434268441052
$ws.Range("Q8").Value = 1.273423344307334
$ws.Range("R8").Value = 11.460810098766
$ws.Range("S8").Value = 0.007732121424641977
$ws.Range("T8").Value = 0.007732121424641976

$ws.Range("G9").Value = 6.636180666666667
$ws.Range("H9").Value = 19.908542
$ws.Range("I9").Value = 0.1078935651261332
$ws.Range("J9").Value = 0.1078935651261332
$ws.Range("O9").Value = 0.4192095090855653
$ws.Range("P9").Value = 0.4192095090855652
$ws.Range("Q9").Value = 7.449048648586889
$ws.Range("R9").Value = 67.041437837282
$ws.Range("S9").Value = 0.04523000847001778
$ws.Range("T9").Value = 0.04523000847001777

$ws.Range("G10").Value = 6.636180666666667
$ws.Range("H10").Value = 19.908542
$ws.Range("I10").Value = 0.1078935651261332
$ws.Range("J10").Value = 0.1078935651261332
$ws.Range("M10").Value = 1.363254333333333
$ws.Range("N10").Value = 4.089763
$ws.Range("O10").Value = 0.5091261482300243
$ws.Range("P10").Value = 0.5091261482300243
$ws.Range("Q10").Value = 9.046802050616222
$ws.Range("R10").Value = 81.421218455546
$ws.Range("S10").Value = 0.05493143523147349
$ws.Range("T10").Value = 0.05493143523147349
